# Update "horarios-141-2026-01-23.xlsx" with the new scrape batch
# (Última actualización: 02:21:47 -> 02:48:52), per commit:
#   "📊 Horarios actualizados Línea 141 - 813"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:48:52"
$ws1.Range("A3").Value = "Total filas: 9"

$ws1.Range("A8").Value = "02:48:52"
$ws1.Range("B8").Value = "02:57"
$ws1.Range("C8").Value = "215_ALUAR"
$ws1.Range("D8").Value = 9
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = "02:21:47"
$ws1.Range("B9").Value = "02:58"
$ws1.Range("C9").Value = "215_ALUAR"
$ws1.Range("D9").Value = 37
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = "01:55:51"
$ws1.Range("B10").Value = "03:12"
$ws1.Range("C10").Value = "215_ALUAR"
$ws1.Range("D10").Value = 77
$ws1.Range("E10").Value = "LP1912"

$ws1.Range("A11").Value = "02:48:52"
$ws1.Range("B11").Value = "03:48"
$ws1.Range("C11").Value = "14_ABASTO"
$ws1.Range("D11").Value = 60
$ws1.Range("E11").Value = "LP1912"

$ws1.Range("A12").Value = "02:21:47"
$ws1.Range("B12").Value = "03:56"
$ws1.Range("C12").Value = "14_ABASTO"
$ws1.Range("D12").Value = 95
$ws1.Range("E12").Value = "LP1912"

$ws1.Range("A13").Value = "02:48:52"
$ws1.Range("B13").Value = "04:01"
$ws1.Range("C13").Value = "81_EL PELIGRO"
$ws1.Range("D13").Value = 73
$ws1.Range("E13").Value = "LP1912"

$ws1.Range("A14").Value = "02:48:52"
$ws1.Range("B14").Value = "04:45"
$ws1.Range("C14").Value = "215A_EL PATO"
$ws1.Range("D14").Value = 117
$ws1.Range("E14").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:48:52"
$ws2.Range("A3").Value = "Total filas: 5"

$ws2.Range("A7").Value = "02:48:52"
$ws2.Range("B7").Value = "02:57"
$ws2.Range("C7").Value = "215_ALUAR"
$ws2.Range("D7").Value = 9
$ws2.Range("E7").Value = "LP1912"

$ws2.Range("A8").Value = "02:21:47"
$ws2.Range("B8").Value = "02:58"
$ws2.Range("C8").Value = "215_ALUAR"
$ws2.Range("D8").Value = 37
$ws2.Range("E8").Value = "LP1912"

$ws2.Range("A9").Value = "01:55:51"
$ws2.Range("B9").Value = "03:12"
$ws2.Range("C9").Value = "215_ALUAR"
$ws2.Range("D9").Value = 77
$ws2.Range("E9").Value = "LP1912"

$ws2.Range("A10").Value = "02:48:52"
$ws2.Range("B10").Value = "04:45"
$ws2.Range("C10").Value = "215A_EL PATO"
$ws2.Range("D10").Value = 117
$ws2.Range("E10").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "6203-6173" — only the "last updated" stamp changes
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 02:48:52"
